$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AO2").Value = 0.03370463745452637
$ws.Range("AP2").Value = -1.7794182328976764
$ws.Range("AO3").Value = 0.033604948995243814
$ws.Range("AP3").Value = 0.8089580161154538
$ws.Range("AO4").Value = 0.034033929934209936
$ws.Range("AP4").Value = 3.3549206958132918
$ws.Range("AO5").Value = 0.034988685198183114
$ws.Range("AP5").Value = 5.745274804981541
$ws.Range("AO6").Value = 0.03646646811566845
$ws.Range("AP6").Value = 7.8896135067488435
$ws.Range("AO7").Value = 0.038464682808313846
$ws.Range("AP7").Value = 9.729635429083386
$ws.Range("AO8").Value = 0.04098088586923354
$ws.Range("AP8").Value = 11.240574825404423
$ws.Range("AO9").Value = 0.04401278841538578
$ws.Range("AP9").Value = 12.426312079030442
$ws.Range("AO10").Value = 0.04755825779785022
$ws.Range("AP10").Value = 13.311198199347896
$ws.Range("AO11").Value = 0.05161531924692535
$ws.Range("AP11").Value = 13.931470699114774
$ws.Range("AO12").Value = 0.05616968963123061
$ws.Range("AP12").Value = 14.327313433194021
$ws.Range("AO13").Value = 0.06114649767023751
$ws.Range("AP13").Value = 14.538828870104407
$ws.Range("AO14").Value = 0.06643006170762508
$ws.Range("AP14").Value = 14.608734752417028
$ws.Range("AO15").Value = 0.07187591055863321
$ws.Range("AP15").Value = 14.578110101245008
$ws.Range("AO16").Value = 0.07731404540004483
$ws.Range("AP16").Value = 14.4830232734512
$ws.Range("AO17").Value = 0.08255340558005234
$ws.Range("AP17").Value = 14.353325022655248
$ws.Range("AO18").Value = 0.08738753035488336
$ws.Range("AP18").Value = 14.212781970451127
$ws.Range("AO19").Value = 0.09160141016152135
$ws.Range("AP19").Value = 14.079873974111907
$ws.Range("AO20").Value = 0.09497952267809454
$ws.Range("AP20").Value = 13.968801025248055
$ws.Range("AO21").Value = 0.09731505060210924
$ws.Range("AP21").Value = 13.89042915054178
$ws.Range("AO22").Value = 0.09842027978838902
$ws.Range("AP22").Value = 13.8530150363786
